# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.301.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "'2.080.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'328.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.4320"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'0.08827"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").Value = "'1.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "'24.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "'2.087.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "'6.721"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'7.660"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "'95.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "'18.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'6.298"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "'30.343.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'12.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "'2.303"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'2.329.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("D27").Value = "'22.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'2.586"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.48%  "
$ws.Range("D29").Value = "'162.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'131.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'1.189"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +20.24%  "
$ws.Range("D34").Value = "'6.173"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "'3.835"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'9.911"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.53%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "'0.06694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "'5.448"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "'12.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("D41").Value = "'0.2261"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").Value = "'0.6805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "'1.243"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'14.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'0.6358"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").Value = "'2.201"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "'3.610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'1.253"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'1.186"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.86%  "
$ws.Range("D51").Value = "'81.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.46%  "
